# Generate Report for Handoff
# - Update the status text from "Handed back: in sync with en-US" to
#   "Ready for handoff" and refresh the associated timestamps.
# - Remove the now-obsolete e9983f9d-... row (row 3) from every sheet,
#   including its hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-16-17 10:16:32"

# Rebuild hyperlinks, dropping the row-3 (e9983f9d...) entry.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e05a09084c5afe10f02e545ad7498e8ff9c2b54/e2e/8fd85048-309c-480e-ab6e-2db41430b6fe.md", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.md")

$ws.Rows(3).Delete()

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-17 10:16:26"

# Rebuild hyperlinks, dropping the row-3 (e9983f9d...) entries.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e05a09084c5afe10f02e545ad7498e8ff9c2b54/e2e/8fd85048-309c-480e-ab6e-2db41430b6fe.md", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e05a09084c5afe10f02e545ad7498e8ff9c2b54/e2e/8fd85048-309c-480e-ab6e-2db41430b6fe.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e8664a98b49ebad87ab7079f5dbd1032e7adbf57/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.zh-cn.xlf", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8c0d274094dcb4ed7fe0a0fd774038b981ad9f5a/e2e/8fd85048-309c-480e-ab6e-2db41430b6fe.md", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0f69475bb99d53d06b5f9491d47fdb798e3c664c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.zh-cn.xlf", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.zh-cn.xlf")

$ws.Rows(3).Delete()

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-17 10:16:32"

# Rebuild hyperlinks, dropping the row-3 (e9983f9d...) entries.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e05a09084c5afe10f02e545ad7498e8ff9c2b54/e2e/8fd85048-309c-480e-ab6e-2db41430b6fe.md", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e05a09084c5afe10f02e545ad7498e8ff9c2b54/e2e/8fd85048-309c-480e-ab6e-2db41430b6fe.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91d6e18df2504acc4ffbfbf9abe128ab26e7a3cc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.de-de.xlf", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d1cb4939b284d8d1a95a66087befb0a570643b66/e2e/8fd85048-309c-480e-ab6e-2db41430b6fe.md", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7579608626ac89d7391e3c43c5f605b82d0ddc8b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.de-de.xlf", "", "", "8fd85048-309c-480e-ab6e-2db41430b6fe.8b5f56a8795213a897f206566c362204bc2938a4.de-de.xlf")

$ws.Rows(3).Delete()
